$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the Genre for "The Richest Man In Babylon" (row 53):
# "Self-Developmet" -> "Self-Development"
$ws.Range("H53").Value = "Self-Development"

# Add the two new books to the end of the list (rows 57 and 58), copying the
# existing formatting from the last data row (56) so the new rows match the
# rest of the table.
$ws.Range("E56:H56").Copy()
$ws.Range("E57:H58").PasteSpecial(-4122)
$ws.Rows.Item(57).RowHeight = 21
$ws.Rows.Item(58).RowHeight = 21

$ws.Range("E57").Value = 53
$ws.Range("F57").Value = "80/20 Principle"
$ws.Range("G57").Value = "Richard Koch"
$ws.Range("H57").Value = "Self-Development"

$ws.Range("E58").Value = 54
$ws.Range("F58").Value = "Why we Sleep"
$ws.Range("G58").Value = "Mathew Walker"
$ws.Range("H58").Value = "Health/Science"

# Update the sheet's selection to match the state after entering the new data.
$ws.Range("L57").Select()
